# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells receive values that would otherwise be auto-parsed as numbers
# (losing trailing zeros / precision); force text format first, matching
# how the source data is stored (plain text strings).
$textCells = @("D4","D5","D8","D9","D11","D13","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.243.56'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '1.883.45'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '314.02'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("D8").Value = '0.3900'
$ws.Range("E8").Value = '  +2.46%  '
$ws.Range("D9").Value = '0.08364'
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("E10").Value = '  +1.13%  '
$ws.Range("D11").Value = '41.57'
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '20.74'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").Value = '1.886.73'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '7.292'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '0.00001107'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '91.43'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").Value = '0.06664'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '17.83'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '6.080'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").Value = '28.286.29'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").Value = '11.17'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '2.277'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '2.094.07'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").Value = '2.515'
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").Value = '159.02'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").Value = '20.66'
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("D30").Value = '125.54'
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").Value = '0.1066'
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D33").Value = '5.878'
$ws.Range("E33").Value = '  +5.26%  '
$ws.Range("D34").Value = '3.593'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '9.770'
$ws.Range("E35").Value = '  +1.41%  '
$ws.Range("D36").Value = '0.02456'
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("D37").Value = '0.06587'
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("D38").Value = '0.2194'
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("E40").Value = '  +2.35%  '
$ws.Range("D41").Value = '5.030'
$ws.Range("E41").Value = '  +3.61%  '
$ws.Range("D42").Value = '1.232'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '11.32'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("D44").Value = '0.6135'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").Value = '13.15'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").Value = '1.291'
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").Value = '3.681'
$ws.Range("D48").Value = '2.014'
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("D49").Value = '1.236'
$ws.Range("E49").Value = '  +2.20%  '
$ws.Range("D50").Value = '121.52'
$ws.Range("D51").Value = '78.99'
$ws.Range("E51").Value = '  -0.87%  '
